# Update cryptos list: refresh Price (D) and Volume(1h) (E) values,
# and fix a handful of scrambled Coin/Link rows (B/C) to match the refreshed ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'29.672.04"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  +0.90%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.928.56"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  +0.84%  "
$c.Style = "Normal"
$c = $ws.Range("D4")
$c.Value = "'1.012"
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  +0.44%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'338.10"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  +4.01%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.009"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  +0.15%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.4821"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.4094"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  +0.72%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.08133"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'1.009"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  -1.11%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'23.53"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  +0.21%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'1.967.56"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  +1.86%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'6.050"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  +0.00%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'7.242"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  +0.45%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'90.66"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -0.42%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  +0.54%  "
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'1.012"
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +0.38%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'0.00001031"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -0.62%  "
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'17.69"
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'1.008"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  +0.05%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'29.693.90"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  +0.92%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'5.607"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -0.36%  "
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'11.84"
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  +0.21%  "
$c.Style = "Normal"
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D24")
$c.Value = "'2.175"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -0.18%  "
$c.Style = "Normal"
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D25")
$c.Value = "'2.124.55"
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -1.22%  "
$c.Style = "Normal"
$ws.Range("B26").Value = "InternetComputer(DFINITY)"
$ws.Range("C26").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D26")
$c.Value = "'6.568"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.22%  "
$c.Style = "Normal"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D27")
$c.Value = "'156.96"
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  +0.73%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'20.00"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -0.24%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'2.082"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -1.13%  "
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'120.78"
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  +0.42%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'1.006"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -1.44%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.09628"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +0.71%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'5.530"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -1.27%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'1.402"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +2.58%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'3.534"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -0.43%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.06564"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +7.48%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.02273"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.57%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'1.197"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  +1.69%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.5941"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -0.62%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'10.71"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -0.98%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'7.904"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  -1.95%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'0.1841"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -0.29%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'2.464"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  +1.99%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'1.271"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  -0.63%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'12.29"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -1.09%  "
$c.Style = "Normal"
$c = $ws.Range("D46")
$c.Value = "'0.07470"
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -1.91%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'0.5538"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -0.64%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'1.975"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  +1.10%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'116.50"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -1.02%  "
$c.Style = "Normal"
$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D50")
$c.Value = "'2.418"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D51")
$c.Value = "'72.17"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -0.02%  "
$c.Style = "Normal"
